# "Generate Report for Handback"
#
# Refreshes the timestamps recorded on the handback status report:
#   - Overview!G2            "Latest HO Xliff Generate Date"       (shared with de-de!H2)
#   - de-de!H2                "Correspond Handoff Datetime"         (shared with Overview!G2)
#   - zh-cn!H2                "Correspond Handoff Datetime"
#   - zh-cn!K2                "Correspond Handback DateTime"
#   - de-de!K2                "Correspond Handback DateTime"
#
# These cells hold plain text (not real Excel date serials), so they are
# written as text the same way the original report generator produced them.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Sheets.Item("Overview")
$wsZhCn     = $wb.Sheets.Item("zh-cn")
$wsDeDe     = $wb.Sheets.Item("de-de")

# Overview / de-de share the same "handoff generated" timestamp for the
# 1656f302-....md file - update both occurrences together.
$wsOverview.Range("G2").Value = "2016-08-31 15:25:35"
$wsDeDe.Range("H2").Value = "2016-08-31 15:25:35"

# zh-cn handoff / handback datetimes for the same file.
$wsZhCn.Range("H2").Value = "2016-08-31 15:25:30"
$wsZhCn.Range("K2").Value = "2016-08-31 15:25:50"

# de-de handback datetime for the same file.
$wsDeDe.Range("K2").Value = "2016-08-31 15:25:57"
